$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F ("想去人数" / want-to-go count)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 34
$ws.Range("F3").Value = 61
$ws.Range("F4").Value = 26
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 143
$ws.Range("F8").Value = 228
$ws.Range("F9").Value = 7055
$ws.Range("F10").Value = 182
$ws.Range("F11").Value = 311
$ws.Range("F12").Value = 5279
$ws.Range("F13").Value = 66
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 1087
$ws.Range("F16").Value = 379
$ws.Range("F17").Value = 382
$ws.Range("F18").Value = 533
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 200
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 102
$ws.Range("F25").Value = 0
$ws.Range("F27").Value = 1828
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 75
$ws.Range("F34").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 5077
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F43").Value = 164
$ws.Range("F45").Value = 1060
$ws.Range("F46").Value = 957
$ws.Range("F49").Value = 0

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 92
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 7

# Sheet "全部类型" (all types) - note: sheet 3 "本地生活" has no data rows, so no changes there
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 1120
$ws.Range("F8").Value = 8686
$ws.Range("F10").Value = 228
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 311
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 5919
$ws.Range("F21").Value = 379
$ws.Range("F22").Value = 382
$ws.Range("F23").Value = 533
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 1554
$ws.Range("F36").Value = 75
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 1990
$ws.Range("F40").Value = 292
$ws.Range("F41").Value = 1318
$ws.Range("F42").Value = 5077
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 606
$ws.Range("F45").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 957
$ws.Range("F50").Value = 1327
$ws.Range("F51").Value = 0
